$wb = $excel.ActiveWorkbook

# weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.35992376256339
$ws.Range("C2").Value = 0.309719915328088
$ws.Range("B3").Value = -0.153144231507661
$ws.Range("C3").Value = 0.150161873039823

# lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.52166620795309
$ws.Range("C2").Value = 0.333572687350795
$ws.Range("B3").Value = -0.773932916603465
$ws.Range("C3").Value = 0.135974053392742

# llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.9482535140164
$ws.Range("C2").Value = 0.320334372530899
$ws.Range("B3").Value = 0.261024068593365
$ws.Range("C3").Value = 0.189629506963794

# gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.2271083115886
$ws.Range("C2").Value = 0.34245155985288
$ws.Range("B3").Value = -0.0377432470081011
$ws.Range("C3").Value = 0.030834685943283

# exp - unchanged, no edits required

# weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0959264259508379
$ws.Range("B2").Value = -0.0218092681903502
$ws.Range("A3").Value = -0.0218092681903502
$ws.Range("B3").Value = 0.0225485881148278

# lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.111270737746431
$ws.Range("B2").Value = -0.0303399210140507
$ws.Range("A3").Value = -0.0303399210140507
$ws.Range("B3").Value = 0.0184889431960524

# llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.102614110224765
$ws.Range("B2").Value = 0.00923537545251059
$ws.Range("A3").Value = 0.00923537545251059
$ws.Range("B3").Value = 0.0359593499113315

# gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.117273070845671
$ws.Range("B2").Value = -0.00652736343257616
$ws.Range("A3").Value = -0.00652736343257616
$ws.Range("B3").Value = 0.000950777857220891

# exp cov - unchanged, no edits required

$wb.Save()
